$wb = $excel.ActiveWorkbook

# Insert the new "AISG Abstract" sheet immediately before "Sheet1", matching the
# target tab order: Generate Cost Function, Sheet2, AISG Time Results,
# AISG Abstract, Sheet1, Experimentation.
$beforeSheet = $wb.Worksheets.Item("Sheet1")
$abstractSheet = $wb.Worksheets.Add($beforeSheet)
$abstractSheet.Name = "AISG Abstract"

$data = @(
    '---------------------------------------------------------------------------------------------                         │ 176572 ipfs       20   0 3634M  484M 16600 S  0.0  2.0  3h44:11 /nix/store/7vkhz7qzqb0b47vcqdh20ih7a3k6xndp-ipfs-0.10Hi: 3000000000000000, Lo: 1, Resolution: 1, Max Time: 5:00:00                                                         │ 259946 ipfs       20   0 3634M  484M 16600 S  0.0  2.0  2h40:42 /nix/store/7vkhz7qzqb0b47vcqdh20ih7a3k6xndp-ipfs-0.10Timeout: 5:00:00, Update Time: 0:00:00.010000, SAT Update Time: 0:00:00.010000                                        │F1Help   F2Setup  F3Search F4Filter F5Tree  F6SortByF7Nice - F8Nice + F9Kill   F10Quit                                Stride discount: 0.5, Stutter discount: 0.75                                                                          │───────────────────────────────────────────────────────────────────────────────────────────────────────────────',
    'Reserved finger: pinky                                                                                                │[jtsoundy@hopper:~/Projects/h-generator/solver/results]$ cat abstract/G_12/5H/config_G_16.5H.abstract.from_G12',
    '---------------------------------------------------------------------------------------------                         │---------------------------------------------------------------------------------------------',
    'N-Grams: 16, Setup Time: 0:00:00.862518, Current Time: 2022-04-28 11:40:49.802445                                     │Hi: 3000000000000000, Lo: 1, Resolution: 1, Max Time: 5:00:00',
    '---------------------------------------------------------------------------------------------                         │Timeout: 5:00:00, Update Time: 0:00:00.010000, SAT Update Time: 0:00:00.010000',
    'Cost Constraint         - Actual Cost             - Result  - Time:This Run  - Time:All Runs                          │Stride discount: 0.5, Stutter discount: 0.75',
    '2,999,999,940,000,000   - 2,567,017,604,161,970   - sat     - 0:00:00.099857 - 0:00:00.099876                         │Reserved finger: pinky',
    '2,567,017,574,161,970   - 2,565,539,991,891,666   - sat     - 0:00:00.036432 - 0:00:00.149505                         │---------------------------------------------------------------------------------------------',
    '2,565,539,961,891,666   - 2,536,623,251,176,254   - sat     - 0:00:00.033910 - 0:00:00.196434                         │Traceback (most recent call last):',
    '2,536,623,221,176,254   - 2,524,705,147,482,716   - sat     - 0:00:00.048200 - 0:00:00.257679                         │  File "/home/jtsoundy/Projects/h-generator/solver/twiddler.py", line 16, in <module>',
    '2,524,705,117,482,716   - 2,518,552,934,789,454   - sat     - 0:00:00.080416 - 0:00:00.351112                         │    b = lib.problem_def(s, n, p)',
    '2,518,552,904,789,454   - 2,494,074,086,794,507   - sat     - 0:00:00.041012 - 0:00:00.405154                         │  File "/home/jtsoundy/Projects/h-generator/solver/lib/buttons.py", line 201, in problem_def',
    '2,494,074,056,794,507   - 2,492,485,230,597,880   - sat     - 0:00:00.902574 - 0:00:01.320765                         │    custom_forbid(G, G_D, s, n)',
    '2,492,485,200,597,880   - 2,483,154,541,211,617   - sat     - 0:00:00.596194 - 0:00:01.930061                         │  File "/home/jtsoundy/Projects/h-generator/solver/lib/buttons.py", line 88, in custom_forbid',
    '2,483,154,511,211,617   - 2,475,531,340,337,118   - sat     - 0:00:00.119272 - 0:00:02.062313                         │    s.add( And( Not(G_D[i] == 4095), Not(G_D[i] == 4088), Not(G_D[i] == 4039), Not(G_D[i] == 3647), Not(G_D[i] == 511)',
    '2,475,531,310,337,118   - 2,438,392,905,369,262   - sat     - 0:00:00.812263 - 0:00:02.887734                         │))',
    '2,438,392,875,369,262   - 2,431,742,038,750,648   - sat     - 0:00:00.414454 - 0:00:03.315331                         │NameError: name ''i'' is not defined',
    '2,431,742,008,750,648   - 2,423,865,935,831,325   - sat     - 0:00:02.499938 - 0:00:05.828451                         │',
    '2,423,865,905,831,325   - 2,423,136,032,844,502   - sat     - 0:00:10.283326 - 0:00:16.124833                         │[jtsoundy@hopper:~/Projects/h-generator/solver/results]$ git fetch',
    '2,423,136,002,844,502   - 2,421,169,188,631,603   - sat     - 0:00:00.410129 - 0:00:16.548178                         │remote: Enumerating objects: 9, done.',
    '2,421,169,158,631,603   - 2,419,852,092,499,990   - sat     - 0:00:01.754067 - 0:00:18.315391                         │remote: Counting objects: 100% (9/9), done.',
    '2,419,852,062,499,990   - 2,419,306,664,045,359   - sat     - 0:00:20.014858 - 0:00:38.343469                         │remote: Compressing objects: 100% (1/1), done.',
    '2,419,306,634,045,359   - 2,410,406,934,676,342   - sat     - 0:00:00.063929 - 0:00:38.420667                         │remote: Total 5 (delta 4), reused 5 (delta 4), pack-reused 0',
    '2,410,406,904,676,342   - 2,410,406,934,676,342   - unsat   - 0:00:23.711567 - 0:01:02.145407                         │Unpacking objects: 100% (5/5), 424 bytes | 212.00 KiB/s, done.',
    '2,410,406,934,676,341   - 2,410,406,934,676,342   - unsat   - 0:00:23.195888 - 0:01:25.359454                         │From https://github.com/HiDefender/h-generator',
    '---------------------------------------------------------------------------------------------                         │   769967b..94023a6  z3-twiddler-model -> origin/z3-twiddler-model',
    'Sat: 2410406934676342, Unknown: 0, Unsat: 2410406934676341.5                                                          │',
    'Total Time: 0:01:26.247111                                                                                            │[jtsoundy@hopper:~/Projects/h-generator/solver/results]$ git reset --hard origin/z3-twiddler-model',
    '---------------------------------------------------------------------------------------------'
)

$values = New-Object 'object[,]' 29,1
for ($i = 0; $i -lt 29; $i++) {
    $values[$i, 0] = $data[$i]
}

$abstractSheet.Range("A1:A29").Value = $values
$abstractSheet.Range("A1:A29").Select()
